$d = $word.ActiveDocument

# 1. Update the intro sentence: remove the comma after "gender", and
#    replace ", and calculates" with ". It calculates"
$d.Content.Find.Execute(
    "gender, and physical activity factor, and calculates",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "gender and physical activity factor. It calculates", 2
) | Out-Null

Write-Output "done"
